$wb = $excel.ActiveWorkbook

$newListFormula = '"1 US-en,2 IN-hi,3 IN-kn,4 es,5 EU-Ge,6 EU-fr,7 Chinese,8 Russian"'

# 1) Update the two data validation list rules on sheet "1 dashboard" and
#    "2 signup" so that the allowed language list also includes the two
#    newly added languages (Chinese and Russian). Using Validation.Modify
#    keeps the original sqref / rule ordering intact and only changes the
#    formula, exactly like the diff shows.
$ws1 = $wb.Worksheets.Item("1 dashboard")
$ws1.Range("C10:C500").Validation.Modify(3, 1, 1, $newListFormula)
$ws1.Range("C2:C500").Validation.Modify(3, 1, 1, $newListFormula)

$ws2 = $wb.Worksheets.Item("2 signup")
$ws2.Range("C10:C500").Validation.Modify(3, 1, 1, $newListFormula)
$ws2.Range("C2:C500").Validation.Modify(3, 1, 1, $newListFormula)

# 2) Add the new worksheet "3 homePage" after the existing sheets. We start
#    from a copy of "1 dashboard" so that the header row (Label_name /
#    Label_value / Language_id), its styling, the column widths and the
#    page setup all match the other sheets exactly, then trim it down to
#    just the header row and refresh the validation lists.
$lastIndex = $wb.Worksheets.Count
$ws1.Copy($null, $wb.Worksheets.Item($lastIndex))
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "3 homePage"

# Remove the sample data rows (2-11), keeping only the header row.
$ws3.Range("A2:A11").EntireRow.Delete()

# Recreate the data validation rules with the updated language list on the
# new sheet (the copy brought over the old formula1 text).
$ws3.Cells.Validation.Delete()

$ws3.Range("C10:C500").Validation.Add(3, 1, 1, $newListFormula)
$ws3.Range("C10:C500").Validation.ShowInput = $false
$ws3.Range("C10:C500").Validation.ShowError = $false

$ws3.Range("C2:C500").Validation.Add(3, 1, 1, $newListFormula)
$ws3.Range("C2:C500").Validation.ShowInput = $false
$ws3.Range("C2:C500").Validation.ShowError = $false


# Restore the original active sheet/selection so the only differences in
# workbook.xml / the existing sheets are the ones described above.
$ws1.Activate()

Write-Host "Added sheet '3 homePage' and refreshed language validation lists."
